# Update the "想去人数" (interested-count) numbers in column F
# on both the "展览" (Exhibition) and "全部类型" (All Types) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$exhibitUpdates = @{
    3  = 5498
    7  = 648
    8  = 626
    9  = 1072
    11 = 1528
    12 = 4970
    13 = 451
    14 = 222
    18 = 4273
    20 = 1141
    22 = 52
    24 = 50
    25 = 154
    29 = 338
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型": row -> new F value
$allUpdates = @{
    4  = 5498
    8  = 648
    9  = 626
    10 = 1072
    12 = 1528
    13 = 4970
    14 = 451
    15 = 222
    19 = 4273
    21 = 1141
    23 = 52
    25 = 50
    26 = 154
    30 = 338
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
